# Scheduled-runner update: refresh currentAveragePrice / derived profit
# columns (H, I, J, K, L, M, N) for a set of leve rows across the ALC,
# ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets, matching the latest
# Universalis price pull. A handful of rows had their computed profit
# collapse to 0 (no market data) so several H/I/J/K/L values become 0
# and the associated M/N cells are cleared entirely.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2191.889
$ws.Range("I100").Value = 1879.5
$ws.Range("J100").Value = 2816.6667
$ws.Range("K100").Value = 1879.5
$ws.Range("L100").Value = 2816.6667
$ws.Range("M100").Value = -1338.5
$ws.Range("N100").Value = -3898.6667

$ws.Range("H113").Value = 2353.5557
$ws.Range("I113").Value = 2061.875
$ws.Range("J113").Value = 2586.9
$ws.Range("K113").Value = 2061.875
$ws.Range("L113").Value = 2586.9
$ws.Range("M113").Value = 1192.125
$ws.Range("N113").Value = -9094.9

$ws.Range("H137").Value = 3339.8235
$ws.Range("I137").Value = 3491.4167
$ws.Range("J137").Value = 2976
$ws.Range("K137").Value = 10474.2501
$ws.Range("L137").Value = 8928
$ws.Range("M137").Value = -7924.250100000001
$ws.Range("N137").Value = -14028

$ws.Range("H138").Value = 2254.747
$ws.Range("I138").Value = 1645.258
$ws.Range("J138").Value = 2592.1428
$ws.Range("K138").Value = 4935.774
$ws.Range("L138").Value = 7776.428400000001
$ws.Range("M138").Value = 204.2259999999997
$ws.Range("N138").Value = -18056.4284

$ws.Range("H141").Value = 4571.485
$ws.Range("I141").Value = 1472.3334
$ws.Range("J141").Value = 7154.1113
$ws.Range("K141").Value = 4417.0002
$ws.Range("L141").Value = 21462.3339
$ws.Range("M141").Value = 762.9997999999996
$ws.Range("N141").Value = -31822.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1235599.9
$ws.Range("I32").Value = 1478449.5
$ws.Range("K32").Value = 1478449.5
$ws.Range("M32").Value = -1478162.5

$ws.Range("H61").Value = 2551.5789
$ws.Range("I61").Value = 2082.6365
$ws.Range("J61").Value = 3196.375
$ws.Range("K61").Value = 2082.6365
$ws.Range("L61").Value = 3196.375
$ws.Range("M61").Value = -1870.6365
$ws.Range("N61").Value = -3620.375

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H97").Value = 926.85
$ws.Range("I97").Value = 832.05884
$ws.Range("J97").Value = 1464
$ws.Range("K97").Value = 832.05884
$ws.Range("L97").Value = 1464
$ws.Range("M97").Value = -336.05884
$ws.Range("N97").Value = -2456

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H102").Value = 2672.8572
$ws.Range("I102").Value = 2242
$ws.Range("K102").Value = 2242
$ws.Range("M102").Value = -620

$ws.Range("H132").Value = 3082.61
$ws.Range("I132").Value = 2394.4792
$ws.Range("K132").Value = 7183.437600000001
$ws.Range("M132").Value = -4653.437600000001

$ws.Range("H136").Value = 2551.5789
$ws.Range("I136").Value = 2082.6365
$ws.Range("J136").Value = 3196.375
$ws.Range("K136").Value = 6247.9095
$ws.Range("L136").Value = 9589.125
$ws.Range("M136").Value = -3697.9095
$ws.Range("N136").Value = -14689.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1731.2307
$ws.Range("I94").Value = 1353.7
$ws.Range("J94").Value = 2989.6667
$ws.Range("K94").Value = 1353.7
$ws.Range("L94").Value = 2989.6667
$ws.Range("M94").Value = -902.7
$ws.Range("N94").Value = -3891.6667

$ws.Range("H134").Value = 6911.636
$ws.Range("I134").Value = 10000
$ws.Range("K134").Value = 30000
$ws.Range("M134").Value = -27465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 924.75
$ws.Range("I16").Value = 950
$ws.Range("J16").Value = 899.5
$ws.Range("K16").Value = 950
$ws.Range("L16").Value = 899.5
$ws.Range("M16").Value = -663
$ws.Range("N16").Value = -1473.5

$ws.Range("H31").Value = 3334.8987
$ws.Range("I31").Value = 1079.0698
$ws.Range("J31").Value = 6029.3613
$ws.Range("K31").Value = 1079.0698
$ws.Range("L31").Value = 6029.3613
$ws.Range("M31").Value = -784.0698
$ws.Range("N31").Value = -6619.3613

$ws.Range("H34").Value = 3334.8987
$ws.Range("I34").Value = 1079.0698
$ws.Range("J34").Value = 6029.3613
$ws.Range("K34").Value = 1079.0698
$ws.Range("L34").Value = 6029.3613
$ws.Range("M34").Value = -877.0698
$ws.Range("N34").Value = -6433.3613

$ws.Range("H107").Value = 891.3889
$ws.Range("I107").Value = 745
$ws.Range("J107").Value = 933.2143
$ws.Range("K107").Value = 745
$ws.Range("L107").Value = 933.2143
$ws.Range("M107").Value = 1175
$ws.Range("N107").Value = -4773.2143

$ws.Range("H113").Value = 924.75
$ws.Range("I113").Value = 950
$ws.Range("J113").Value = 899.5
$ws.Range("K113").Value = 950
$ws.Range("L113").Value = 899.5
$ws.Range("M113").Value = 1220
$ws.Range("N113").Value = -5239.5

$ws.Range("H122").Value = 1769.1305
$ws.Range("J122").Value = 1798.4762
$ws.Range("L122").Value = 5395.4286
$ws.Range("N122").Value = -10295.4286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1966.3043
$ws.Range("J5").Value = 2129.75
$ws.Range("L5").Value = 6389.25
$ws.Range("N5").Value = -6613.25

$ws.Range("H68").Value = 6647
$ws.Range("J68").Value = 1080.7693
$ws.Range("L68").Value = 3242.3079
$ws.Range("N68").Value = -4864.3079

$ws.Range("H71").Value = 6647
$ws.Range("J71").Value = 1080.7693
$ws.Range("L71").Value = 9726.923699999999
$ws.Range("N71").Value = -17838.9237

$ws.Range("H76").Value = 3999.3547
$ws.Range("J76").Value = 3999.6667
$ws.Range("L76").Value = 11999.0001
$ws.Range("N76").Value = -12765.0001

$ws.Range("H79").Value = 3999.3547
$ws.Range("J79").Value = 3999.6667
$ws.Range("L79").Value = 11999.0001
$ws.Range("N79").Value = -14651.0001

$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()

$ws.Range("H109").Value = 2694.4092
$ws.Range("I109").Value = 571.5454999999999
$ws.Range("J109").Value = 4817.273
$ws.Range("K109").Value = 1714.6365
$ws.Range("L109").Value = 14451.819
$ws.Range("M109").Value = -674.6364999999998
$ws.Range("N109").Value = -16531.819

$ws.Range("H121").Value = 1144.9736
$ws.Range("I121").Value = 401.4
$ws.Range("J121").Value = 1410.5358
$ws.Range("K121").Value = 1204.2
$ws.Range("L121").Value = 4231.607400000001
$ws.Range("M121").Value = 105.8000000000002
$ws.Range("N121").Value = -6851.607400000001

$ws.Range("H122").Value = 3800.25
$ws.Range("I122").Value = 459.75
$ws.Range("J122").Value = 9367.75
$ws.Range("K122").Value = 4137.75
$ws.Range("L122").Value = 84309.75
$ws.Range("M122").Value = -1687.75
$ws.Range("N122").Value = -89209.75

$ws.Range("H135").Value = 1966.3043
$ws.Range("J135").Value = 2129.75
$ws.Range("L135").Value = 19167.75
$ws.Range("N135").Value = -24237.75

$ws.Range("H136").Value = 3847.0952
$ws.Range("I136").Value = 1348.1666
$ws.Range("K136").Value = 4044.4998
$ws.Range("M136").Value = 1055.5002

$ws.Range("H141").Value = 6141.963
$ws.Range("I141").Value = 3062.2
$ws.Range("J141").Value = 9991.666999999999
$ws.Range("K141").Value = 9186.599999999999
$ws.Range("L141").Value = 29975.001
$ws.Range("M141").Value = -4006.599999999999
$ws.Range("N141").Value = -40335.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 795
$ws.Range("I97").Value = 795
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 795
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -299
$ws.Range("N97").ClearContents()

$ws.Range("H107").Value = 424
$ws.Range("J107").Value = 400
$ws.Range("L107").Value = 400
$ws.Range("N107").Value = -4240

$ws.Range("H122").Value = 2358
$ws.Range("I122").Value = 2358
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7074
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4624
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 1582.6
$ws.Range("I126").Value = 1304.3334
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 3913.0002
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -1443.0002
$ws.Range("N126").Value = -10940

$ws.Range("H132").Value = 3386
$ws.Range("I132").Value = 2629.7646
$ws.Range("J132").Value = 7671.3335
$ws.Range("K132").Value = 7889.293799999999
$ws.Range("L132").Value = 23014.0005
$ws.Range("M132").Value = -5359.293799999999
$ws.Range("N132").Value = -28074.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4066583.2
$ws.Range("I136").Value = 1570.8529
$ws.Range("K136").Value = 4712.5587
$ws.Range("M136").Value = -2162.5587

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 805.1818
$ws.Range("I107").Value = 837.6667
$ws.Range("J107").Value = 735.5714
$ws.Range("K107").Value = 2513.0001
$ws.Range("L107").Value = 2206.7142
$ws.Range("M107").Value = -593.0001000000002
$ws.Range("N107").Value = -6046.7142

$ws.Range("H122").Value = 2134.3416
$ws.Range("I122").Value = 2071.4
$ws.Range("J122").Value = 2501.5
$ws.Range("K122").Value = 6214.200000000001
$ws.Range("L122").Value = 7504.5
$ws.Range("M122").Value = -3764.200000000001
$ws.Range("N122").Value = -12404.5
